$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Professional summary paragraph: neutralize voter-impact language
# ---------------------------------------------------------------------
$pSummary = $d.Paragraphs(4)
$pSummary.Range.Find.Execute(
    "affecting all Black and Asian-American voters", $true, $false, $false,
    $false, $false, $true, 1, $false, "affecting 50M voters", 2)

# ---------------------------------------------------------------------
# 2) Siege Analytics bullet: split the run so "50M" is its own
#    bold / colored run, matching the other stat call-outs in the bullet
# ---------------------------------------------------------------------
$pBullet = $d.Paragraphs(10)
$pBullet.Range.Find.Execute(
    "affecting all Black and Asian-American voters,", $true, $false, $false,
    $false, $false, $true, 1, $false, "affecting 50M voters,", 2)

$numberRange = $pBullet.Range.Duplicate
$numberRange.Find.Execute("50M", $true, $false, $false, $false, $false,
                           $true, 1, $false, "", 0)
$numberRange.Font.Bold = 1
$numberRange.Font.Color = 5258796

# ---------------------------------------------------------------------
# 3) Move the "Software Engineer - Mautinoa Technologies" block from
#    after Salsa Labs up to right after the Siege Analytics bullets
#    (i.e. immediately before "Senior Analyst - Myers Research")
# ---------------------------------------------------------------------
$i = 1
$mautinoaStart = 0
$mautinoaEnd = 0
$seniorAnalystIdx = 0
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Software Engineer - Mautinoa Technologies*") {
        $mautinoaStart = $i
        $mautinoaEnd = $i + 4
    }
    if ($t -like "Senior Analyst - Myers Research*") {
        $seniorAnalystIdx = $i
    }
    $i = $i + 1
}

$srcFirst = $d.Paragraphs($mautinoaStart)
$srcLast = $d.Paragraphs($mautinoaEnd)
$moveRange = $d.Range($srcFirst.Range.Start, $srcLast.Range.End)
$moveRange.Cut()

$destPara = $d.Paragraphs($seniorAnalystIdx)
$destRange = $d.Range($destPara.Range.Start, $destPara.Range.Start)
$destRange.Paste()

# The paste drops the Heading3 paragraph style from the moved heading;
# restore it.
$d.Paragraphs($seniorAnalystIdx).Style = "Heading3"

# ---------------------------------------------------------------------
# 4) Key Projects "Impact" line: neutralize voter-impact language
# ---------------------------------------------------------------------
$impactFound = $false
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Impact: Corrected demographic data affecting*") {
        $p.Range.Find.Execute(
            "affecting all Black and Asian-American voters", $true, $false,
            $false, $false, $false, $true, 1, $false,
            "affecting 50M voters nationwide", 2)
        $impactFound = $true
    }
}

Write-Output "done"
